# Applies the cryptos-list price/volume(1h) refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text even when it looks numeric (e.g. "215.35"),
# matching the source data which stores these as text strings, then drop the
# temporary "@" (text) number format back to Normal so no stray cell style sticks.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "25.810.21"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.635.62"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue "D5" "215.35"
$ws.Range("E5").Value = "  -0.24%  "
Set-TextValue "D6" "0.5053"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.20%  "
Set-TextValue "D8" "0.2579"
$ws.Range("E8").Value = "  +0.18%  "
Set-TextValue "D9" "0.06416"
$ws.Range("E9").Value = "  +1.03%  "
Set-TextValue "D10" "20.29"
$ws.Range("E10").Value = "  +3.84%  "
Set-TextValue "D11" "0.07798"
$ws.Range("E11").Value = "  +0.60%  "
Set-TextValue "D12" "4.291"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.635.05"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.861.06"
$ws.Range("E14").Value = "  +0.14%  "
Set-TextValue "D15" "0.5632"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "0.0₅7634"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "25.827.57"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  +0.08%  "
Set-TextValue "D20" "194.34"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -1.76%  "
Set-TextValue "D22" "9.920"
$ws.Range("E22").Value = "  +0.21%  "
Set-TextValue "D23" "6.100"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("E24").Value = "  +0.08%  "
Set-TextValue "D25" "1.786"
$ws.Range("E25").Value = "  -6.00%  "
Set-TextValue "D26" "140.35"
$ws.Range("E26").Value = "  -1.43%  "
Set-TextValue "D27" "0.1255"
$ws.Range("E27").Value = "  +1.42%  "
Set-TextValue "D28" "6.821"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  -0.96%  "
Set-TextValue "D31" "0.04909"
$ws.Range("E31").Value = "  +0.62%  "
Set-TextValue "D32" "3.314"
$ws.Range("E32").Value = "  +2.13%  "
Set-TextValue "D33" "3.243"
$ws.Range("E33").Value = "  +1.81%  "
Set-TextValue "D34" "1.576"
$ws.Range("E34").Value = "  +2.47%  "
Set-TextValue "D35" "2.379"
$ws.Range("E35").Value = "  +0.03%  "
Set-TextValue "D36" "0.9058"
$ws.Range("E36").Value = "  +0.65%  "
Set-TextValue "D37" "2.570"
$ws.Range("E37").Value = "  +0.41%  "
Set-TextValue "D38" "0.5541"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "1.124.28"
$ws.Range("E40").Value = "  +0.40%  "
Set-TextValue "D42" "5.517"
$ws.Range("E42").Value = "  -1.03%  "
Set-TextValue "D43" "0.8009"
$ws.Range("E43").Value = "  -0.34%  "
Set-TextValue "D44" "98.20"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").Value = "1.771.03"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  -5.71%  "
Set-TextValue "D47" "55.54"
$ws.Range("E47").Value = "  +1.42%  "
Set-TextValue "D48" "0.4265"
$ws.Range("E48").Value = "  -4.15%  "
Set-TextValue "D49" "7.728"
$ws.Range("E49").Value = "  +3.11%  "
Set-TextValue "D50" "0.05018"
$ws.Range("E50").Value = "  -2.56%  "
Set-TextValue "D51" "1.002"
$ws.Range("E51").Value = "  +0.03%  "
